$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")
$ws.Activate()

# Row 8 = item 7: "Implementar Swagger OpenAPI" -> "Implementar Swagger OpenAPI e testes de integração"
$ws.Range("B8").Value = "Implementar Swagger OpenAPI e testes de integração"

# The longer text now wraps onto a second line, so the row grows taller
$ws.Rows.Item(8).RowHeight = 30.75

# Selection moved to J11 on re-save
$ws.Range("J11").Select()
